# March 24 update 3
# Adds three new columns (M: renewd, N: PlanID, O: iteration) to the sheet,
# mirroring the header style used by the other header cells, and fills in
# the corresponding values for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), matching the bold/bordered header style ---
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)

# --- Fill new columns for each existing data row (rows 2-17) ---
$lastRow = 17
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"
    $ws.Cells.Item($r, 14).Value = 20131419
    $ws.Cells.Item($r, 15).Value = 13
}
